# proAutomobileInsurance.xlsx - "new radiobutton control" demo-app update
#
# Adds the "FillPageInsurantData" step to the Tabelle1 record/dialog table
# (new shared string + C2 cell), nudges the B/C column widths to fit the
# new content, and leaves the selection on H16 (matching the author's
# last-saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New row-2 entry in column C: the InsurantData smoke-test step.
$ws.Range("C2").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageInsurantData"

# Re-fit columns B and C now that C holds the long step name.
# (Target widths are 60.6640625 / 60.77734375 "chars"; ColumnWidth rounds to
# the nearest whole pixel on write, so these inputs are chosen to land on the
# closest reachable pixel-quantized width to each target.)
$ws.Columns.Item(2).ColumnWidth = 59.8333333333333
$ws.Columns.Item(3).ColumnWidth = 60

# Leave the cursor where the author left it when saving.
$ws.Range("H16").Select() | Out-Null
